# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (columns H-N) for specific
# Leve rows across the ALC, ARM, BSM, CRP, GSM, LTW and WVR sheets.
# Values mirror a scheduled data-refresh run; some rows gain or lose the
# optional LeveProfitNQ (M) / LeveProfitHQ (N) cells depending on whether a
# profit is computable for NQ vs HQ pricing this cycle.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3671.2144
$ws.Range("I74").Value = 3448.8333
$ws.Range("J74").Value = 3838
$ws.Range("K74").Value = 3448.8333
$ws.Range("L74").Value = 3838
$ws.Range("M74").Value = -2512.8333
$ws.Range("N74").Value = -5710
# Row 77
$ws.Range("H77").Value = 3671.2144
$ws.Range("I77").Value = 3448.8333
$ws.Range("J77").Value = 3838
$ws.Range("K77").Value = 17244.1665
$ws.Range("L77").Value = 19190
$ws.Range("M77").Value = -12564.1665
$ws.Range("N77").Value = -28550
# Row 106
$ws.Range("H106").Value = 22225242
$ws.Range("I106").Value = 27780552
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 27780552
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -27779921
$ws.Range("N106").Value = -5262
# Row 124
$ws.Range("H124").Value = 28052.5
$ws.Range("I124").Value = 10210
$ws.Range("J124").Value = 34000
$ws.Range("K124").Value = 10210
$ws.Range("L124").Value = 34000
$ws.Range("M124").Value = -5300
$ws.Range("N124").Value = -43820
# Row 129
$ws.Range("H129").Value = 1212.0667
$ws.Range("J129").Value = 1277.4286
$ws.Range("L129").Value = 3832.2858
$ws.Range("N129").Value = -13832.2858
# Row 130
$ws.Range("H130").Value = 98888
$ws.Range("J130").Value = 98888
$ws.Range("L130").Value = 98888
$ws.Range("N130").Value = -108928
# Row 135
$ws.Range("H135").Value = 1032.7858
$ws.Range("I135").Value = 1035.3077
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 9317.7693
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -6782.7693
$ws.Range("N135").Value = -14070
# Row 138
$ws.Range("H138").Value = 6726130.5
$ws.Range("I138").Value = 2979291.5
$ws.Range("J138").Value = 8200296.5
$ws.Range("K138").Value = 8937874.5
$ws.Range("L138").Value = 24600889.5
$ws.Range("M138").Value = -8932734.5
$ws.Range("N138").Value = -24611169.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 253852.75
$ws.Range("I2").Value = 253852.75
$ws.Range("K2").Value = 253852.75
$ws.Range("M2").Value = -253739.75
# Row 32
$ws.Range("H32").Value = 36984.574
$ws.Range("I32").Value = 8729.361000000001
$ws.Range("J32").Value = 129456.18
$ws.Range("K32").Value = 8729.361000000001
$ws.Range("L32").Value = 129456.18
$ws.Range("M32").Value = -8442.361000000001
$ws.Range("N32").Value = -130030.18
# Row 97
$ws.Range("H97").Value = 5561.75
$ws.Range("I97").Value = 6337.353
$ws.Range("K97").Value = 6337.353
$ws.Range("M97").Value = -5841.353
# Row 116
$ws.Range("H116").Value = 253852.75
$ws.Range("I116").Value = 253852.75
$ws.Range("K116").Value = 253852.75
$ws.Range("M116").Value = -251558.75
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 253852.75
$ws.Range("I3").Value = 253852.75
$ws.Range("K3").Value = 253852.75
$ws.Range("M3").Value = -253738.75
# Row 20
$ws.Range("H20").Value = 4000
$ws.Range("I20").Value = 4000
$ws.Range("J20").Value = 4000
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 4000
$ws.Range("M20").Value = -3753
$ws.Range("N20").Value = -4494
# Row 105
$ws.Range("H105").Value = 3611.7144
$ws.Range("I105").Value = 3654.5833
$ws.Range("J105").Value = 3518.182
$ws.Range("K105").Value = 3654.5833
$ws.Range("L105").Value = 3518.182
$ws.Range("M105").Value = -1907.5833
$ws.Range("N105").Value = -7012.182
# Row 132
$ws.Range("H132").Value = 45271.11
$ws.Range("J132").Value = 45271.11
$ws.Range("L132").Value = 45271.11
$ws.Range("N132").Value = -55391.11
# Row 134
$ws.Range("H134").Value = 1962.5641
$ws.Range("I134").Value = 1875.3889
$ws.Range("J134").Value = 3008.6667
$ws.Range("K134").Value = 5626.1667
$ws.Range("L134").Value = 9026.000100000001
$ws.Range("M134").Value = -3091.1667
$ws.Range("N134").Value = -14096.0001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6357.9
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 6357.9
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 6357.9
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -6947.9
# Row 34
$ws.Range("H34").Value = 6357.9
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 6357.9
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6357.9
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -6761.9
# Row 105
$ws.Range("H105").Value = 884.3125
$ws.Range("I105").Value = 751.92
$ws.Range("J105").Value = 1357.1428
$ws.Range("K105").Value = 751.92
$ws.Range("L105").Value = 1357.1428
$ws.Range("M105").Value = 995.08
$ws.Range("N105").Value = -4851.1428
# Row 138
$ws.Range("H138").Value = 57695
$ws.Range("J138").Value = 57695
$ws.Range("L138").Value = 57695
$ws.Range("N138").Value = -67975

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2918.2856
$ws.Range("I126").Value = 2771.2
$ws.Range("K126").Value = 8313.599999999999
$ws.Range("M126").Value = -5843.599999999999
# Row 132
$ws.Range("H132").Value = 3805.24
$ws.Range("I132").Value = 3788.2666
$ws.Range("K132").Value = 11364.7998
$ws.Range("M132").Value = -8834.799800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 94
$ws.Range("H94").Value = 59999.5
$ws.Range("J94").Value = 59999.5
$ws.Range("L94").Value = 59999.5
$ws.Range("N94").Value = -61351.5

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -3746
# Row 101
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
# Row 125
$ws.Range("H125").Value = 29860
$ws.Range("J125").Value = 29860
$ws.Range("L125").Value = 29860
$ws.Range("N125").Value = -39700
# Row 132
$ws.Range("H132").Value = 2498.8071
$ws.Range("I132").Value = 2464.8538
$ws.Range("J132").Value = 2585.8125
$ws.Range("K132").Value = 7394.5614
$ws.Range("L132").Value = 7757.4375
$ws.Range("M132").Value = -4864.5614
$ws.Range("N132").Value = -12817.4375
# Row 136
$ws.Range("H136").Value = 2610.8572
$ws.Range("I136").Value = 723.4828
$ws.Range("J136").Value = 11733.167
$ws.Range("K136").Value = 2170.4484
$ws.Range("L136").Value = 35199.501
$ws.Range("M136").Value = 379.5515999999998
$ws.Range("N136").Value = -40299.501
